$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 265.47058
$ws.Range("I12").Value = 264.85715
$ws.Range("K12").Value = 264.85715
$ws.Range("M12").Value = -94.85714999999999

# Row 29
$ws.Range("H29").Value = 16050
$ws.Range("I29").Value = 3750
$ws.Range("J29").Value = 22200
$ws.Range("K29").Value = 11250
$ws.Range("L29").Value = 66600
$ws.Range("M29").Value = -10969
$ws.Range("N29").Value = -67162

# Row 41
$ws.Range("H41").Value = 500
$ws.Range("I41").Value = 500
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 500
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -60
$ws.Range("N41").ClearContents()

# Row 69
$ws.Range("H69").Value = 4000
$ws.Range("I69").Value = 4000
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 12000
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -11126
$ws.Range("N69").ClearContents()

# Row 72
$ws.Range("H72").Value = 4000
$ws.Range("I72").Value = 4000
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 36000
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -31632
$ws.Range("N72").ClearContents()

# Row 76
$ws.Range("H76").Value = 9638.454
$ws.Range("I76").Value = 9894.444
$ws.Range("J76").Value = 8486.5
$ws.Range("K76").Value = 9894.444
$ws.Range("L76").Value = 8486.5
$ws.Range("M76").Value = -9579.444
$ws.Range("N76").Value = -9116.5

# Row 79
$ws.Range("H79").Value = 9638.454
$ws.Range("I79").Value = 9894.444
$ws.Range("J79").Value = 8486.5
$ws.Range("K79").Value = 9894.444
$ws.Range("L79").Value = 8486.5
$ws.Range("M79").Value = -8802.444
$ws.Range("N79").Value = -10670.5

# Row 132
$ws.Range("H132").Value = 4670.9355
$ws.Range("I132").Value = 5216.423
$ws.Range("K132").Value = 15649.269
$ws.Range("M132").Value = -13119.269

# Row 141
$ws.Range("H141").Value = 1478.3939
$ws.Range("I141").Value = 1249.0358
$ws.Range("J141").Value = 2762.8
$ws.Range("K141").Value = 3747.1074
$ws.Range("L141").Value = 8288.400000000001
$ws.Range("M141").Value = 1432.8926
$ws.Range("N141").Value = -18648.4

$ws = $wb.Worksheets.Item("ARM")
# Row 122
$ws.Range("H122").Value = 2147.3
$ws.Range("J122").Value = 2582.5
$ws.Range("L122").Value = 7747.5
$ws.Range("N122").Value = -12647.5

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1714.4286
$ws.Range("I94").Value = 1423.8667
$ws.Range("J94").Value = 2049.6924
$ws.Range("K94").Value = 1423.8667
$ws.Range("L94").Value = 2049.6924
$ws.Range("M94").Value = -972.8667
$ws.Range("N94").Value = -2951.6924

$ws = $wb.Worksheets.Item("CRP")
# Row 105
$ws.Range("H105").Value = 10807.5
$ws.Range("I105").Value = 2642.111
$ws.Range("J105").Value = 35303.668
$ws.Range("K105").Value = 2642.111
$ws.Range("L105").Value = 35303.668
$ws.Range("M105").Value = -895.1109999999999
$ws.Range("N105").Value = -38797.668

# Row 141
$ws.Range("H141").Value = 87129.89
$ws.Range("J141").Value = 93646.125
$ws.Range("L141").Value = 93646.125
$ws.Range("N141").Value = -104006.125

$ws = $wb.Worksheets.Item("CUL")
# Row 55
$ws.Range("H55").Value = 8409696
$ws.Range("I55").Value = 2250150
$ws.Range("J55").Value = 13337332
$ws.Range("K55").Value = 6750450
$ws.Range("L55").Value = 40011996
$ws.Range("M55").Value = -6750273
$ws.Range("N55").Value = -40012350

# Row 81
$ws.Range("H81").Value = 1251435.8
$ws.Range("I81").Value = 5000745
$ws.Range("J81").Value = 1666
$ws.Range("K81").Value = 15002235
$ws.Range("L81").Value = 4998
$ws.Range("M81").Value = -15001112
$ws.Range("N81").Value = -7244

# Row 84
$ws.Range("H84").Value = 1251435.8
$ws.Range("I84").Value = 5000745
$ws.Range("J84").Value = 1666
$ws.Range("K84").Value = 45006705
$ws.Range("L84").Value = 14994
$ws.Range("M84").Value = -45001089
$ws.Range("N84").Value = -26226

# Row 86
$ws.Range("H86").Value = 1451.5385
$ws.Range("I86").Value = 984.5
$ws.Range("K86").Value = 2953.5
$ws.Range("M86").Value = -1767.5

# Row 89
$ws.Range("H89").Value = 1451.5385
$ws.Range("I89").Value = 984.5
$ws.Range("K89").Value = 8860.5
$ws.Range("M89").Value = -2932.5

# Row 122
$ws.Range("H122").Value = 654.65216
$ws.Range("I122").Value = 300
$ws.Range("J122").Value = 670.7727
$ws.Range("K122").Value = 2700
$ws.Range("L122").Value = 6036.954299999999
$ws.Range("M122").Value = -250
$ws.Range("N122").Value = -10936.9543

# Row 130
$ws.Range("H130").Value = 2574
$ws.Range("J130").Value = 3088.6667
$ws.Range("L130").Value = 9266.000100000001
$ws.Range("N130").Value = -19306.0001

$ws = $wb.Worksheets.Item("GSM")
# Row 93
$ws.Range("H93").Value = 33331
$ws.Range("I93").Value = 25000
$ws.Range("J93").Value = 35711.285
$ws.Range("K93").Value = 25000
$ws.Range("L93").Value = 35711.285
$ws.Range("N93").Value = -39455.285
$ws.Range("M93").Value = -23128

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1553.2222
$ws.Range("I16").Value = 1396.85
$ws.Range("K16").Value = 1396.85
$ws.Range("M16").Value = -1226.85

# Row 68
$ws.Range("H68").Value = 2031.6666
$ws.Range("I68").Value = 1496.6666
$ws.Range("J68").Value = 2566.6667
$ws.Range("K68").Value = 1496.6666
$ws.Range("L68").Value = 2566.6667
$ws.Range("N68").Value = -4064.6667
$ws.Range("M68").Value = -747.6666

# Row 71
$ws.Range("H71").Value = 2031.6666
$ws.Range("I71").Value = 1496.6666
$ws.Range("J71").Value = 2566.6667
$ws.Range("K71").Value = 7483.333000000001
$ws.Range("L71").Value = 12833.3335
$ws.Range("N71").Value = -20321.3335
$ws.Range("M71").Value = -3739.333000000001

# Row 93
$ws.Range("H93").Value = 2305.7273
$ws.Range("I93").Value = 2320.3447
$ws.Range("K93").Value = 2320.3447
$ws.Range("M93").Value = -1072.3447

# Row 100
$ws.Range("H100").Value = 3928.8667
$ws.Range("I100").Value = 3402.5833
$ws.Range("K100").Value = 3402.5833
$ws.Range("M100").Value = -2861.5833

# Row 122
$ws.Range("H122").Value = 4333.9165
$ws.Range("I122").Value = 4117
$ws.Range("J122").Value = 4637.6
$ws.Range("K122").Value = 12351
$ws.Range("L122").Value = 13912.8
$ws.Range("M122").Value = -9901
$ws.Range("N122").Value = -18812.8

$ws = $wb.Worksheets.Item("WVR")
# Row 75
$ws.Range("H75").Value = 49993.6
$ws.Range("J75").Value = 49992.25
$ws.Range("L75").Value = 49992.25
$ws.Range("N75").Value = -51864.25

# Row 78
$ws.Range("H78").Value = 49993.6
$ws.Range("J78").Value = 49992.25
$ws.Range("L78").Value = 149976.75
$ws.Range("N78").Value = -159336.75

# Row 136
$ws.Range("H136").Value = 10712.826
$ws.Range("I136").Value = 9446.708000000001
$ws.Range("K136").Value = 28340.124
$ws.Range("M136").Value = -25790.124
